$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose numeric-looking text would otherwise be auto-converted to a
# number by Excel (losing a significant trailing zero, e.g. "580.50" -> 580.5).
# Force them to Text format first so the literal string is preserved, matching
# the source workbook which stores these as text values.
$textCells = @("D5", "D8", "D30", "D44")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.771.32"
$ws.Range("E2").Value = "  +3.14%  "
$ws.Range("D3").Value = "2.529.43"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "580.50"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "152.81"
$ws.Range("E6").Value = "  +3.53%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.540"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").Value = "2.530.74"
$ws.Range("E9").Value = "  +2.71%  "
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "29.23"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "2.989.89"
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("D17").Value = "64.570.29"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "2.533.78"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "7.94"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "11.01"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("E21").Value = "  +3.57%  "
$ws.Range("D22").Value = "330.23"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "10.13"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "65.82"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").Value = "645.07"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  +6.69%  "
$ws.Range("D29").Value = "2.641.59"
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("D30").Value = "1.50"
$ws.Range("E30").Value = "  +4.47%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").Value = "8.07"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "1.57"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("D37").Value = "4.86"
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("D39").Value = "155.17"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").Value = "18.98"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "1.81"
$ws.Range("E43").Value = "  +4.04%  "
$ws.Range("D44").Value = "163.10"
$ws.Range("E44").Value = "  +7.20%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "0.0₆0302"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").Value = "21.52"
$ws.Range("E49").Value = "  +5.03%  "
$ws.Range("D50").Value = "0.623"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("E51").Value = "  +1.38%  "
